$wb = $excel.ActiveWorkbook

$wsCore = $wb.Worksheets.Item(1)   # Core basis vectors relative to ...
$wsAlpha = $wb.Worksheets.Item(2)  # Alpha, Beta to Normal vector re...
$wsMap = $wb.Worksheets.Item(3)    # Normal vector relative to Map g...
$wsDip = $wb.Worksheets.Item(4)    # Normal vector to Dip, Azimuth

# --- Step 1: copy the {d,r,s} basis vectors (relative to x,y,z) computed on
# the "Core basis vectors" sheet into the "Normal vector relative to Map g"
# sheet's input cells.
$wsCore.Activate()
$wsCore.Range("P5:X5").Select()
$wsCore.Range("P5:X5").Copy()
$wsMap.Activate()
$wsMap.Range("A5:I5").Select()
$wsMap.Range("A5:I5").PasteSpecial(-4163)

# --- Step 2: copy the normal vector n (relative to {d,r,s}) computed on the
# "Alpha, Beta to Normal vector" sheet into the same "Normal vector relative
# to Map g" sheet's remaining input cells.
$wsAlpha.Activate()
$wsAlpha.Range("G5:I5").Select()
$wsAlpha.Range("G5:I5").Copy()
$wsMap.Activate()
$wsMap.Range("J5:L5").Select()
$wsMap.Range("J5:L5").PasteSpecial(-4163)

# --- Step 3: copy the resulting normal vector n (relative to {x,y,z}) into
# the "Normal vector to Dip, Azimuth" sheet's input cells.
$wsMap.Activate()
$wsMap.Range("M5:O5").Select()
$wsMap.Range("M5:O5").Copy()
$wsDip.Activate()
$wsDip.Range("A5:C5").Select()
$wsDip.Range("A5:C5").PasteSpecial(-4163)

# --- Final selections left on each sheet (source anchor cell remains the
# active cell on the sheets we copied from; the paste anchor remains active
# on the destination sheet).
$wsCore.Range("P5").Select()
$wsAlpha.Range("G5").Select()
$wsMap.Range("M5").Select()
$wsDip.Range("A5").Select()

# "Normal vector to Dip, Azimuth" ends up the active sheet/tab.
$wsDip.Activate()
